# Append new pupae/fly data rows (118-146) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows to append, starting at row 118.
# Each entry: vial, treatment (5=conditioned, 6=unconditioned), time_hours, males, females
$data = @(
    @(1,  "conditioned",   332, 3, 3),
    @(1,  "unconditioned", 332, 0, 0),
    @(2,  "conditioned",   332, 5, 6),
    @(2,  "unconditioned", 332, 0, 0),
    @(3,  "conditioned",   332, 5, 8),
    @(3,  "unconditioned", 332, 0, 0),
    @(4,  "conditioned",   332, 5, 3),
    @(4,  "unconditioned", 332, 4, 1),
    @(5,  "conditioned",   332, 0, 0),
    @(5,  "unconditioned", 332, 5, 3),
    @(6,  "conditioned",   332, 2, 1),
    @(6,  "unconditioned", 332, 4, 1),
    @(7,  "conditioned",   332, 0, 0),
    @(7,  "unconditioned", 332, 4, 4),
    @(8,  "conditioned",   332, 3, 5),
    @(8,  "unconditioned", 332, 3, 2),
    @(9,  "conditioned",   332, 0, 1),
    @(9,  "unconditioned", 332, 4, 8),
    @(10, "conditioned",   332, 4, 0),
    @(10, "unconditioned", 332, 8, 2),
    @(11, "conditioned",   332, 3, 1),
    @(11, "unconditioned", 332, 3, 1),
    @(12, "conditioned",   332, 1, 0),
    @(12, "unconditioned", 332, 4, 4),
    @(13, "conditioned",   332, 0, 0),
    @(13, "unconditioned", 332, 1, 0),
    @(14, "conditioned",   332, 1, 1),
    @(14, "unconditioned", 332, 1, 0)
)

$row = 118
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Final row 146 only has vial, treatment, time_hours (no males/females counts yet).
$ws.Cells.Item(146, 1).Value = 15
$ws.Cells.Item(146, 2).Value = "unconditioned"
$ws.Cells.Item(146, 3).Value = 332

# Update the view to match the saved state in the diff (scrolled/selected/zoomed).
$ws.Range("D146").Select()
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 82
